# Natmi following Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt1"
$ws.Cells.Item(2, 3).Value = "Fzd8"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01948966666666667
$ws.Cells.Item(2, 8).Value = 0.058469
$ws.Cells.Item(2, 9).Value = 0.0709606244933031
$ws.Cells.Item(2, 10).Value = 0.0709606244933031
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.062569
$ws.Cells.Item(2, 14).Value = 9.187707
$ws.Cells.Item(2, 15).Value = 0.2460077391286943
$ws.Cells.Item(2, 16).Value = 0.2690593624267
$ws.Cells.Item(2, 17).Value = 0.05968844895366666
$ws.Cells.Item(2, 18).Value = 0.537196040583
$ws.Cells.Item(2, 19).Value = 0.01745686279875775
$ws.Cells.Item(2, 20).Value = 0.01909262038356861

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt1"
$ws.Cells.Item(3, 3).Value = "Fzd8"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01948966666666667
$ws.Cells.Item(3, 8).Value = 0.058469
$ws.Cells.Item(3, 9).Value = 0.0709606244933031
$ws.Cells.Item(3, 10).Value = 0.0709606244933031
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.147102
$ws.Cells.Item(3, 14).Value = 18.441306
$ws.Cells.Item(3, 15).Value = 0.4937797859292232
$ws.Cells.Item(3, 16).Value = 0.5400483531609875
$ws.Cells.Item(3, 17).Value = 0.119804968946
$ws.Cells.Item(3, 18).Value = 1.078244720514
$ws.Cells.Item(3, 19).Value = 0.0350389219717072
$ws.Cells.Item(3, 20).Value = 0.03832216839688358

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt1"
$ws.Cells.Item(4, 3).Value = "Fzd8"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.01948966666666667
$ws.Cells.Item(4, 8).Value = 0.058469
$ws.Cells.Item(4, 9).Value = 0.0709606244933031
$ws.Cells.Item(4, 10).Value = 0.0709606244933031
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.01759833333333333
$ws.Cells.Item(4, 14).Value = 0.052795
$ws.Cells.Item(4, 15).Value = 0.001413625683459368
$ws.Cells.Item(4, 16).Value = 0.001546086421706485
$ws.Cells.Item(4, 17).Value = 0.0003429856505555555
$ws.Cells.Item(4, 18).Value = 0.003086870855
$ws.Cells.Item(4, 19).Value = 0.0001003117612980491
$ws.Cells.Item(4, 20).Value = 0.0001097112580049086

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt1"
$ws.Cells.Item(5, 3).Value = "Fzd8"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.01948966666666667
$ws.Cells.Item(5, 8).Value = 0.058469
$ws.Cells.Item(5, 9).Value = 0.0709606244933031
$ws.Cells.Item(5, 10).Value = 0.0709606244933031
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.022088
$ws.Cells.Item(5, 14).Value = 0.066264
$ws.Cells.Item(5, 15).Value = 0.001774268250568265
$ws.Cells.Item(5, 16).Value = 0.00194052222081558
$ws.Cells.Item(5, 17).Value = 0.0004304877573333333
$ws.Cells.Item(5, 18).Value = 0.003874389816
$ws.Cells.Item(5, 19).Value = 0.0001259031830789645
$ws.Cells.Item(5, 20).Value = 0.000137700668632205

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Wnt1"
$ws.Cells.Item(6, 3).Value = "Fzd8"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.01948966666666667
$ws.Cells.Item(6, 8).Value = 0.058469
$ws.Cells.Item(6, 9).Value = 0.0709606244933031
$ws.Cells.Item(6, 10).Value = 0.0709606244933031
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.1997185
$ws.Cells.Item(6, 14).Value = 6.399437000000001
$ws.Cells.Item(6, 15).Value = 0.2570245810080548
$ws.Cells.Item(6, 16).Value = 0.1874056757697904
$ws.Cells.Item(6, 17).Value = 0.06236144699216667
$ws.Cells.Item(6, 18).Value = 0.374168681953
$ws.Cells.Item(6, 19).Value = 0.01823862477846114
$ws.Cells.Item(6, 20).Value = 0.01329842378621381

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt1"
$ws.Cells.Item(7, 3).Value = "Fzd8"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.255165
$ws.Cells.Item(7, 8).Value = 0.7654949999999999
$ws.Cells.Item(7, 9).Value = 0.9290393755066968
$ws.Cells.Item(7, 10).Value = 0.9290393755066969
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.062569
$ws.Cells.Item(7, 14).Value = 9.187707
$ws.Cells.Item(7, 15).Value = 0.2460077391286943
$ws.Cells.Item(7, 16).Value = 0.2690593624267
$ws.Cells.Item(7, 17).Value = 0.7814604188849998
$ws.Cells.Item(7, 18).Value = 7.033143769964999
$ws.Cells.Item(7, 19).Value = 0.2285508763299365
$ws.Cells.Item(7, 20).Value = 0.2499667420431314

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt1"
$ws.Cells.Item(8, 3).Value = "Fzd8"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.255165
$ws.Cells.Item(8, 8).Value = 0.7654949999999999
$ws.Cells.Item(8, 9).Value = 0.9290393755066968
$ws.Cells.Item(8, 10).Value = 0.9290393755066969
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 6.147102
$ws.Cells.Item(8, 14).Value = 18.441306
$ws.Cells.Item(8, 15).Value = 0.4937797859292232
$ws.Cells.Item(8, 16).Value = 0.5400483531609875
$ws.Cells.Item(8, 17).Value = 1.56852528183
$ws.Cells.Item(8, 18).Value = 14.11672753647
$ws.Cells.Item(8, 19).Value = 0.458740863957516
$ws.Cells.Item(8, 20).Value = 0.501726184764104

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt1"
$ws.Cells.Item(9, 3).Value = "Fzd8"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.255165
$ws.Cells.Item(9, 8).Value = 0.7654949999999999
$ws.Cells.Item(9, 9).Value = 0.9290393755066968
$ws.Cells.Item(9, 10).Value = 0.9290393755066969
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.01759833333333333
$ws.Cells.Item(9, 14).Value = 0.052795
$ws.Cells.Item(9, 15).Value = 0.001413625683459368
$ws.Cells.Item(9, 16).Value = 0.001546086421706485
$ws.Cells.Item(9, 17).Value = 0.004490478724999999
$ws.Cells.Item(9, 18).Value = 0.040414308525
$ws.Cells.Item(9, 19).Value = 0.001313313922161318
$ws.Cells.Item(9, 20).Value = 0.001436375163701577

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Wnt1"
$ws.Cells.Item(10, 3).Value = "Fzd8"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.255165
$ws.Cells.Item(10, 8).Value = 0.7654949999999999
$ws.Cells.Item(10, 9).Value = 0.9290393755066968
$ws.Cells.Item(10, 10).Value = 0.9290393755066969
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.022088
$ws.Cells.Item(10, 14).Value = 0.066264
$ws.Cells.Item(10, 15).Value = 0.001774268250568265
$ws.Cells.Item(10, 16).Value = 0.00194052222081558
$ws.Cells.Item(10, 17).Value = 0.005636084519999999
$ws.Cells.Item(10, 18).Value = 0.05072476067999999
$ws.Cells.Item(10, 19).Value = 0.0016483650674893
$ws.Cells.Item(10, 20).Value = 0.001802821552183375

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Wnt1"
$ws.Cells.Item(11, 3).Value = "Fzd8"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.255165
$ws.Cells.Item(11, 8).Value = 0.7654949999999999
$ws.Cells.Item(11, 9).Value = 0.9290393755066968
$ws.Cells.Item(11, 10).Value = 0.9290393755066969
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.1997185
$ws.Cells.Item(11, 14).Value = 6.399437000000001
$ws.Cells.Item(11, 15).Value = 0.2570245810080548
$ws.Cells.Item(11, 16).Value = 0.1874056757697904
$ws.Cells.Item(11, 17).Value = 0.8164561710525
$ws.Cells.Item(11, 18).Value = 4.898737026315
$ws.Cells.Item(11, 19).Value = 0.2387859562295936
$ws.Cells.Item(11, 20).Value = 0.1741072519835766

